$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 12005.444
$ws.Range("I41").Value = 929.4
$ws.Range("J41").Value = 25850.5
$ws.Range("K41").Value = 929.4
$ws.Range("L41").Value = 25850.5
$ws.Range("M41").Value = -489.4
$ws.Range("N41").Value = -26730.5

$ws.Range("H76").Value = 3099
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 3099
$ws.Range("K76").Value = 0
$ws.Range("L76").ClearContents()
$ws.Range("M76").Value = 3099
$ws.Range("N76").Value = -3729

$ws.Range("H79").Value = 3099
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 3099
$ws.Range("K79").Value = 0
$ws.Range("L79").ClearContents()
$ws.Range("M79").Value = 3099
$ws.Range("N79").Value = -5283

$ws.Range("H103").Value = 2399.3333
$ws.Range("I103").Value = 1798.8
$ws.Range("K103").Value = 5396.4
$ws.Range("M103").Value = -4810.4

$ws.Range("H135").Value = 15630557
$ws.Range("I135").Value = 21742010
$ws.Range("J135").Value = 12399.444
$ws.Range("K135").Value = 195678090
$ws.Range("L135").Value = 111594.996
$ws.Range("M135").Value = -195675555
$ws.Range("N135").Value = -116664.996

$ws.Range("H137").Value = 5564236.5
$ws.Range("I137").Value = 7417743
$ws.Range("J137").Value = 3717
$ws.Range("K137").Value = 22253229
$ws.Range("L137").Value = 11151
$ws.Range("M137").Value = -22250679
$ws.Range("N137").Value = -16251

$ws.Range("H138").Value = 4257.143
$ws.Range("J138").Value = 4616.15
$ws.Range("L138").Value = 13848.45
$ws.Range("N138").Value = -24128.45

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3729.8
$ws.Range("I61").Value = 3347.3333
$ws.Range("K61").Value = 3347.3333
$ws.Range("M61").Value = -3135.3333

$ws.Range("H74").Value = 1843.2667
$ws.Range("I74").Value = 1884.4286
$ws.Range("J74").Value = 1807.25
$ws.Range("K74").Value = 1884.4286
$ws.Range("L74").Value = 1807.25
$ws.Range("M74").Value = -1010.4286
$ws.Range("N74").Value = -3555.25

$ws.Range("H77").Value = 1843.2667
$ws.Range("I77").Value = 1884.4286
$ws.Range("J77").Value = 1807.25
$ws.Range("K77").Value = 9422.143
$ws.Range("L77").Value = 9036.25
$ws.Range("M77").Value = -5054.143
$ws.Range("N77").Value = -17772.25

$ws.Range("H122").Value = 5129.213
$ws.Range("I122").Value = 4838.92
$ws.Range("J122").Value = 5459.091
$ws.Range("K122").Value = 14516.76
$ws.Range("L122").Value = 16377.273
$ws.Range("M122").Value = -12066.76
$ws.Range("N122").Value = -21277.273

$ws.Range("H132").Value = 8903.348
$ws.Range("I132").Value = 12087.375
$ws.Range("K132").Value = 36262.125
$ws.Range("M132").Value = -33732.125

$ws.Range("H136").Value = 3729.8
$ws.Range("I136").Value = 3347.3333
$ws.Range("K136").Value = 10041.9999
$ws.Range("M136").Value = -7491.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4135.375
$ws.Range("J94").Value = 5188
$ws.Range("L94").Value = 5188
$ws.Range("N94").Value = -6090

$ws.Range("H105").Value = 1618.5385
$ws.Range("I105").Value = 860
$ws.Range("K105").Value = 860
$ws.Range("M105").Value = 887

$ws.Range("H134").Value = 2965.0908
$ws.Range("I134").Value = 3036.7
$ws.Range("J134").Value = 2249
$ws.Range("K134").Value = 9110.099999999999
$ws.Range("L134").Value = 6747
$ws.Range("M134").Value = -6575.099999999999
$ws.Range("N134").Value = -11817

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 358
$ws.Range("I2").Value = 346
$ws.Range("K2").Value = 346
$ws.Range("M2").Value = -233

$ws.Range("H31").Value = 5063.683
$ws.Range("I31").Value = 5829.731
$ws.Range("K31").Value = 5829.731
$ws.Range("M31").Value = -5534.731

$ws.Range("H34").Value = 5063.683
$ws.Range("I34").Value = 5829.731
$ws.Range("K34").Value = 5829.731
$ws.Range("M34").Value = -5627.731

$ws.Range("H112").Value = 20246.8
$ws.Range("J112").Value = 20246.8
$ws.Range("L112").Value = 20246.8
$ws.Range("N112").Value = -23200.8

$ws.Range("H132").Value = 750
$ws.Range("I132").Value = 750
$ws.Range("K132").Value = 2250
$ws.Range("M132").Value = 280

$ws.Range("H134").Value = 2646.4595
$ws.Range("I134").Value = 2987.2415
$ws.Range("J134").Value = 1411.125
$ws.Range("K134").Value = 8961.7245
$ws.Range("L134").Value = 4233.375
$ws.Range("M134").Value = -6426.7245
$ws.Range("N134").Value = -9303.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 75111550
$ws.Range("I4").Value = 48952700
$ws.Range("K4").Value = 146858100
$ws.Range("M4").Value = -146857988

$ws.Range("H57").Value = 202224.5
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 202224.5
$ws.Range("K57").Value = 0
$ws.Range("L57").ClearContents()
$ws.Range("M57").Value = 606673.5
$ws.Range("N57").Value = -607791.5

$ws.Range("H68").Value = 2213.4375
$ws.Range("I68").Value = 1400
$ws.Range("J68").Value = 2267.6667
$ws.Range("K68").Value = 4200
$ws.Range("L68").Value = 6803.000100000001
$ws.Range("M68").Value = -3389
$ws.Range("N68").Value = -8425.000100000001

$ws.Range("H71").Value = 2213.4375
$ws.Range("I71").Value = 1400
$ws.Range("J71").Value = 2267.6667
$ws.Range("K71").Value = 12600
$ws.Range("L71").Value = 20409.0003
$ws.Range("M71").Value = -8544
$ws.Range("N71").Value = -28521.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J57").Value = 15000
$ws.Range("L57").Value = 15000
$ws.Range("N57").Value = -16640

$ws.Range("H80").Value = 76200.2
$ws.Range("I80").Value = 76200.2
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 76200.2
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -75202.2

$ws.Range("H83").Value = 76200.2
$ws.Range("I83").Value = 76200.2
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 381001
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -376009

$ws.Range("H126").Value = 5203.6
$ws.Range("I126").Value = 5254.5
$ws.Range("K126").Value = 15763.5
$ws.Range("M126").Value = -13293.5

$ws.Range("H132").Value = 8243.875
$ws.Range("I132").Value = 7635.9287
$ws.Range("K132").Value = 22907.7861
$ws.Range("M132").Value = -20377.7861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5505.2
$ws.Range("I7").Value = 5075.25
$ws.Range("K7").Value = 5075.25
$ws.Range("M7").Value = -4963.25

$ws.Range("H61").Value = 3639.8572
$ws.Range("I61").Value = 3413.1667
$ws.Range("K61").Value = 3413.1667
$ws.Range("M61").Value = -3211.1667

$ws.Range("H113").Value = 3639.8572
$ws.Range("I113").Value = 3413.1667
$ws.Range("K113").Value = 3413.1667
$ws.Range("M113").Value = -1243.1667

$ws.Range("H123").Value = 30600
$ws.Range("J123").Value = 30600
$ws.Range("L123").Value = 30600
$ws.Range("N123").Value = -40400

$ws.Range("H126").Value = 5505.2
$ws.Range("I126").Value = 5075.25
$ws.Range("K126").Value = 15225.75
$ws.Range("M126").Value = -12755.75

$ws.Range("H132").Value = 5303.3237
$ws.Range("I132").Value = 5487.524
$ws.Range("J132").Value = 5005.769
$ws.Range("K132").Value = 16462.572
$ws.Range("L132").Value = 15017.307
$ws.Range("M132").Value = -13932.572
$ws.Range("N132").Value = -20077.307

$ws.Range("H136").Value = 6112.2856
$ws.Range("I136").Value = 3356.2
$ws.Range("J136").Value = 13002.5
$ws.Range("K136").Value = 10068.6
$ws.Range("L136").Value = 39007.5
$ws.Range("M136").Value = -7518.599999999999
$ws.Range("N136").Value = -44107.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 34900
$ws.Range("I8").Value = 34900
$ws.Range("K8").Value = 34900
$ws.Range("M8").Value = -34760

$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31248

$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96240

$ws.Range("H96").Value = 63300.555
$ws.Range("I96").Value = 90450.164
$ws.Range("K96").Value = 90450.164
$ws.Range("M96").Value = -89077.164

$ws.Range("H107").Value = 382.2857
$ws.Range("I107").Value = 375.41666
$ws.Range("K107").Value = 1126.24998
$ws.Range("M107").Value = 793.7500199999999

$ws.Range("H122").Value = 1475
$ws.Range("I122").Value = 1475
$ws.Range("J122").Value = 1475
$ws.Range("K122").Value = 4425
$ws.Range("L122").Value = 4425
$ws.Range("M122").Value = -1975
$ws.Range("N122").Value = -9325

$ws.Range("H132").Value = 2961.4878
$ws.Range("I132").Value = 3179.423
$ws.Range("K132").Value = 9538.269
$ws.Range("M132").Value = -7008.269

$ws.Range("H136").Value = 3605.1714
$ws.Range("I136").Value = 3564.6287
$ws.Range("K136").Value = 10693.8861
$ws.Range("M136").Value = -8143.8861

$ws.Range("H140").Value = 52333.168
$ws.Range("J140").Value = 52333.168
$ws.Range("L140").Value = 52333.168
$ws.Range("N140").Value = -62693.168
